$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that should be bumped from 45179 to 45180
# for every data row (rows 2 through 482).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 482
}

$ws.Range("C2:C$lastRow").Value = 45180
